$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row of the crypto table has a Price (column D) and a
# Volume(1h) (column E) cell, both stored as plain text. Column D
# values often look numeric ("303.95", "1.002", ...), so Excel
# would silently coerce a plain .Value assignment into a Number.
# Force the cell to Text first (and restore the Normal style right
# after) so the written value round-trips as a string, matching the
# original workbook formatting exactly.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.405.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.629.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3760"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3647"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.93"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08212"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.231"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.529"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001248"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.328"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.629.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.503"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.403.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.136"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.452"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.300"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.811.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.253"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.804"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.041"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02782"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2518"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08744"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07092"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.006"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7032"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.345"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6545"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.292"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.975"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08008"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.46%  "
